$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")

# New row to add at the bottom of the data table
$newRow = 4

$ws.Cells.Item($newRow, 1).Value = "econ_program_unitcost_vaccination"
$ws.Cells.Item($newRow, 2).Value = "yes"
$ws.Cells.Item($newRow, 3).Value = 1
$ws.Cells.Item($newRow, 4).Value = "yes"

# Columns E..BE (5..57) all hold the same numeric value of 82
for ($col = 5; $col -le 57; $col++) {
    $ws.Cells.Item($newRow, $col).Value = 82
}

# Update the column A width to fit the new, longer label
$ws.Columns.Item(1).ColumnWidth = 32.666667

# Move / update the active selection as recorded after the edit
$ws.Range("BD7").Select()
